# #5: cash & deposit done
#
# Rewrites the "存款" (deposits) sheet (sheet4) from its old ad-hoc layout
# (bank / deposit_type / currency / owner / (blank) / total) into the
# standard property-table layout used by the other sheets in this workbook:
#   bank | deposit_type | currency | owner | total | property_category |
#   category | date | legislator_name | legislator_id | source_file | index
# with a proper header row and per-row metadata columns (H:M) that mirror
# the pattern already used on the 汽車/股票/... sheets.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("存款")

# ----------------------------------------------------------------------
# 1. Header row (row 1): real column names instead of the old row mirrored
#    straight off row 2's data.
# ----------------------------------------------------------------------

# Columns B:G already exist with the header style (s=1) from the template;
# just overwrite their values.
$ws.Range("B1").Value2 = "bank"
$ws.Range("C1").Value2 = "deposit_type"
$ws.Range("D1").Value2 = "currency"
$ws.Range("E1").Value2 = "owner"
$ws.Range("F1").Value2 = "total"
$ws.Range("G1").Value2 = "property_category"

# Columns H:M are new - broadcast the header style from B1 across them first
# (Range.Copy onto a multi-cell target replicates formatting to every cell),
# then fill in the real values.
$ws.Range("B1").Copy($ws.Range("H1:M1"))
$ws.Range("H1").Value2 = "category"
$ws.Range("I1").Value2 = "date"
$ws.Range("J1").Value2 = "legislator_name"
$ws.Range("K1").Value2 = "legislator_id"
$ws.Range("L1").Value2 = "source_file"
$ws.Range("M1").Value2 = "index"

# ----------------------------------------------------------------------
# 2. Data rows (2:14): re-populate every column with the corrected data.
# ----------------------------------------------------------------------

# index, bank, deposit_type, currency, owner, total
$rows = @(
    @(51, "台北富邦商業銀行", "活期儲蓄存款", "新臺幣", "丁守中", 1766028),
    @(52, "中華郵政股份有限公司", "活期儲蓄存款", "新臺幣", "丁守中", 302001),
    @(53, "永豐商業銀行", "活期儲蓄存款", "新臺幣", "丁守中", 2795),
    @(54, "台北富邦商業銀行", "活期儲蓄存款", "新臺幣", "溫子苓", 3066404),
    @(55, "中華郵政股份有限公司", "活期存款", "新臺幣", "溫子苓", 27282),
    @(56, "台北富邦商業銀行", "定期存款", "美金", "溫子苓", 29.43),
    @(57, "台北富邦商業銀行", "支票存款", "新臺幣", "溫子苓", 10000),
    @(58, "上海商業儲蓄銀行", "活期儲蓄存款", "新臺幣", "溫子苓", 314213),
    @(59, "國泰世華商業銀行", "活期存款", "新臺幣", "溫子苓", 109441),
    @(60, "聯邦商業銀行", "活期儲蓄存款", "新臺幣", "溫子苓", 5000),
    @(61, "兆豐國際商業銀行", "活期儲蓄存款", "新臺幣", "溫子苓", 92882),
    @(62, "華南商業銀行", "活期儲蓄存款", "新臺幣", "丁守中", 5994),
    @(63, "華南商業銀行", "活期儲蓄存款", "新臺幣", "溫子苓", 121334)
)

# Broadcast the data-row style (s=2, taken from the existing B2) across the
# new H:M columns for every data row in one shot before filling values in.
$ws.Range("B2").Copy($ws.Range("H2:M14"))

for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $i + 2
    $idx = $rows[$i][0]
    $bank = $rows[$i][1]
    $dtype = $rows[$i][2]
    $currency = $rows[$i][3]
    $owner = $rows[$i][4]
    $total = $rows[$i][5]

    $ws.Range("A$r").Value2 = $idx
    $ws.Range("B$r").Value2 = $bank
    $ws.Range("C$r").Value2 = $dtype
    $ws.Range("D$r").Value2 = $currency
    $ws.Range("E$r").Value2 = $owner
    $ws.Range("F$r").Value2 = $total
    $ws.Range("G$r").Value2 = "deposit"
    $ws.Range("H$r").Value2 = "normal"
    $ws.Range("I$r").Value2 = "2013-12-26"
    $ws.Range("J$r").Value2 = "丁守中"
    $ws.Range("K$r").Value2 = 515
    $ws.Range("L$r").Value2 = "tmpc7fb1"
    $ws.Range("M$r").Value2 = $idx
}
